# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sorts / refreshes the "Estado de Cuenta" worker-arrears table (rows
# 16-34, cols B:G on Hoja1): data is now ordered period-ascending (1810,
# 1812, 1901..1905) with each period grouping the three workers, and the
# "Salario Basico" (col G) / some "Valor Mora" (col F) figures were
# updated for the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# row -> Tipo Doc, N Doc, Nombre, Periodo Mora, Valor Mora, Salario Basico
$rows = @{
    16 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1810", 53615, 1000000)
    17 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1812", 53615, 1000000)
    18 = @("CC", "8834880",    "ALCIDES JESUS CARDENAS LOPEZ", "1812", 56472, 1000000)
    19 = @("CC", "1070822062", "DEYMER RAMOS LOPEZ",           "1812", 31249, 781242)
    20 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1901", 53615, 1000000)
    21 = @("CC", "8834880",    "ALCIDES JESUS CARDENAS LOPEZ", "1901", 56472, 1000000)
    22 = @("CC", "1070822062", "DEYMER RAMOS LOPEZ",           "1901", 31249, 781242)
    23 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1902", 53615, 1000000)
    24 = @("CC", "8834880",    "ALCIDES JESUS CARDENAS LOPEZ", "1902", 56472, 1000000)
    25 = @("CC", "1070822062", "DEYMER RAMOS LOPEZ",           "1902", 31249, 781242)
    26 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1903", 53615, 1000000)
    27 = @("CC", "8834880",    "ALCIDES JESUS CARDENAS LOPEZ", "1903", 56472, 1000000)
    28 = @("CC", "1070822062", "DEYMER RAMOS LOPEZ",           "1903", 31249, 781242)
    29 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1904", 53615, 1000000)
    30 = @("CC", "8834880",    "ALCIDES JESUS CARDENAS LOPEZ", "1904", 40000, 1000000)
    31 = @("CC", "1070822062", "DEYMER RAMOS LOPEZ",           "1904", 31249, 781242)
    32 = @("CC", "73476371",   "MARCOS RICARDO GARCES VILLA",  "1905", 26666, 1000000)
    33 = @("CC", "8834880",    "ALCIDES JESUS CARDENAS LOPEZ", "1905", 26666, 1000000)
    34 = @("CC", "1070822062", "DEYMER RAMOS LOPEZ",           "1905", 20833, 781242)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $vals[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $vals[5]   # G - Salario Basico
}
